$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.278.17'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '3.494.33'
$ws.Range("E3").Value = '  -0.29%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '588.62'
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.84'
$ws.Range("E6").Value = '  -1.26%  '

$ws.Range("E8").Value = '  -0.72%  '

$ws.Range("E9").Value = '  +6.56%  '

$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.388'
$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").Value = '4.088.04'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("E13").Value = '  +0.07%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000179'
$ws.Range("E14").Value = '  -1.52%  '

$ws.Range("D15").Value = '3.493.43'
$ws.Range("E15").Value = '  -0.27%  '

$ws.Range("D16").Value = '64.217.35'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("E17").Value = '  -4.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.03'
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.76'
$ws.Range("E19").Value = '  -0.35%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.50'
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '385.24'
$ws.Range("E21").Value = '  -2.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.579'
$ws.Range("E22").Value = '  +1.40%  '

$ws.Range("D23").Value = '3.633.31'
$ws.Range("E23").Value = '  -0.33%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.38'
$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.25'
$ws.Range("E29").Value = '  -3.05%  '

$ws.Range("E30").Value = '  +0.37%  '

$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.99'
$ws.Range("E32").Value = '  -3.76%  '

$ws.Range("E33").Value = '  +3.70%  '

$ws.Range("D34").Value = '3.523.20'
$ws.Range("E34").Value = '  -0.10%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.14'
$ws.Range("E36").Value = '  -1.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.33'
$ws.Range("E37").Value = '  +1.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.86'
$ws.Range("E38").Value = '  -0.89%  '

$ws.Range("E39").Value = '  -2.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.86'
$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0781'
$ws.Range("E41").Value = '  -0.43%  '

$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("E44").Value = '  -0.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.17'
$ws.Range("E45").Value = '  -0.50%  '

$ws.Range("E46").Value = '  -6.41%  '

$ws.Range("E47").Value = '  -1.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.935'
$ws.Range("E48").Value = '  +4.29%  '

$ws.Range("E49").Value = '  -0.30%  '

$ws.Range("D50").Value = '2.390.29'
$ws.Range("E50").Value = '  -3.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0257'
$ws.Range("E51").Value = '  -2.04%  '
